# Alta Cliente y Actualizar Cliente Back-end prontos
#
# This edit renames the "Sprint 1" task labels that were actually used for
# Sprint 3 work (rows 33-35 of Hoja1) so they read "Sprint 3 ..." instead of
# "Sprint 1 ...", and appends a new log row for the "Alta cliente" /
# "Conexión a la BD" back-end work that was just finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fix the mislabeled "Sprint 1" rows -> "Sprint 3" ------------------
# D33: "Sprint 1 - Diagrama de Clases" -> "Sprint 3 - Diagrama de Clases"
$ws.Range("D33").Value = "Sprint 3 - Diagrama de Clases"

# D34 & D35: "Sprint 1 - Back-end" -> "Sprint 3 - Back-end"
$ws.Range("D34").Value = "Sprint 3 - Back-end"
$ws.Range("D35").Value = "Sprint 3 - Back-end"

# --- Append the new work-log row (row 43) -------------------------------
$ws.Range("A43").Value = "Federico Speroni"
$ws.Range("B42").Copy()
$ws.Range("B43").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing date style
$ws.Range("B43").Value = 42870
$ws.Range("C43").Value = 3
$ws.Range("D43").Value = "Sprint 3 - BackEnd"
$ws.Range("E43").Value = "Conexión a la BD, alta cliente"

# --- Update the view state so the new row is visible/selected ----------
$ws.Range("A44").Select()
$excel.ActiveWindow.ScrollRow = 28
